$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Summary" paragraph: merge the many runs that make up the first
#    paragraph into a single run, and change "2009" -> "year 2007" in the
#    "time series analysis ranging from ..." sentence.
# ---------------------------------------------------------------------------
$en_dash = [char]0x2013

$old1 = "In peer-to-peer lending, there are three main stakeholders: borrowers, lenders and the company itself. In my Tableau story I have done exploration on the relationship between these people, what affects borrowers Prosper Score and who defaults the most. First, I have done a time series analysis ranging from 2009 " + $en_dash + " 2014 about the number of loans taken by borrowers, the amount of their loans and how their "
$new1 = "In peer-to-peer lending, there are three main stakeholders: borrowers, lenders and the company itself. In my Tableau story I have done exploration on the relationship between these people, what affects borrowers Prosper Score and who defaults the most. First, I have done a time series analysis ranging from year 2007 " + $en_dash + " 2014 about the number of loans taken by borrowers, the amount of their loans and how their "

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Second paragraph ("Then I have explored the defaults...") - merge all
#    of its runs into one and apply the wording tweaks.
# ---------------------------------------------------------------------------
$lsquo = [char]0x2018
$rsquo = [char]0x2019

$old2 = "Then I have explored the defaults, I found out that the people with " + "`$0 income have highest default rates and most defaulters invest in the loan type " + $lsquo + "D" + $rsquo + ". Breaking down to occupation wise, an interesting pattern was found that the college student group which are enrolled in higher grade studies have more loans, higher borrower and default rate. While this made sense with " + "`$0 income, the sophomore students were the top defaulters and having lower number of loans. Lastly, I have looked at the income and losses on different loan ratings " + $en_dash + " the " + $lsquo + "HR" + $rsquo + " loan rating had the highest loss even though this type of loan is given to most credit-worthy borrowers. However, looking at the net principal returns over the time I noticed loans C&D had highest losses than other loans and are most risky."
$new2 = "Then I have explored the defaulters, reason for defaults, and reason for borrowers to take loan, I found out that the people with " + "`$0 income have highest default rates and most defaulters invest in the loan type " + $lsquo + "D" + $rsquo + ". Breaking down to occupation-wise, an interesting pattern was found that the college student group which are enrolled in higher grade studies have more loans, higher borrower and default rates. While this made sense with " + "`$0 income, the sophomore students were the top defaulters and having lower number of loans. Lastly, I have looked at the incomes and losses on different loan ratings " + $en_dash + " the " + $lsquo + "HR" + $rsquo + " loan rating had the highest loss even though this type of loan is given to most credit-worthy borrowers. However, looking at the net principal returns over the time I noticed loans C&D had highest losses than other loans and are most risky."

$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the empty paragraph that used to sit between the summary text
#    and the "Design Decisions" heading.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $p.Next().Range.Text.StartsWith("Design Decisions")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the "Showing state-wise default rate"
#    bullet up onto the "Design Decisions" heading, and merge the runs
#    around the old bookmark location / around the heading run.
# ---------------------------------------------------------------------------
$headingRange = $d.Content
$headingRange.Find.Execute("Design Decisions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingStart = $headingRange.Start
$bmRange = $d.Range($headingStart, $headingStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 5) Clean up the bullet point text (this also absorbs / removes the old
#    "_GoBack" bookmark that used to live in the middle of this sentence).
# ---------------------------------------------------------------------------
$old3 = " this I created a new numerical variable called Default Rate from Loan Status and mapped with Borrower State, this showed an interesting finding that in states like CA, TX, NY, IL the default rates were quite high with CA having highest defaulters (>700)."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null
